$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = 9
$ws.Range("C12").Value = "2 sessions, preprocessing session 4"
$ws.Range("C12").Select()
